# Edit script for la_mesure_de_prudence.docx
#
# This reproduces two changes:
#   1. In the title paragraph, reorder the "_GoBack" bookmark so that it
#      sits right after the title run (collapsed there) instead of at the
#      very start of the paragraph, and put the "_Toc..." bookmark's end
#      tag right before it (this is what Word naturally does once the
#      insertion point / last-edit marker has moved past the heading text).
#   2. Justify (both-align) every paragraph in the body except the title
#      paragraph and the blank paragraph that immediately follows it.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Move the hidden "_GoBack" bookmark from the start of the heading
#    paragraph to right after the heading text (collapsed position),
#    while keeping it nested so that "_Toc96536214" still closes first.
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titleLen = $titlePara.Range.End - $titlePara.Range.Start - 1   # exclude paragraph mark

if ($d.Bookmarks.Exists("_GoBack")) {
    $goBack = $d.Bookmarks.Item("_GoBack")
    $goBack.Delete()

    $insertionPoint = $d.Range($titleLen, $titleLen)
    $insertionPoint.InsertBefore("X")

    $placeholder = $d.Range($titleLen, $titleLen + 1)
    $d.Bookmarks.Add("_GoBack", $placeholder)
    $placeholder.Text = ""
}

# ---------------------------------------------------------------------
# 2) Justify every paragraph except the title (#1) and the following
#    blank paragraph (#2).
# ---------------------------------------------------------------------
$wdAlignParagraphJustify = 3
$count = $d.Paragraphs.Count
for ($i = 3; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $para.Format.Alignment = $wdAlignParagraphJustify
}
